$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / Changed date) holds the same date serial (45181 = 2023-09-12)
# for every data row (rows 2-171). Update it to 45182 (2023-09-13) for all of them.
$ws.Range("C2:C171").Value = 45182
